$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws4.Name = "Sheet4"

$ws4.Range("C3").Value = "Area"
$ws4.Range("D3").Value = "RL"
$ws4.Range("E3").Value = "WL"
$ws4.Range("F3").Value = "RE"
$ws4.Range("G3").Value = "WE"
$ws4.Range("H3").Value = "LP"

$ws4.Range("B2").Value = "Bigtable"

$ws4.Range("B4").Value = "Area"
$ws4.Range("C4").Value = 3.06
$ws4.Range("D4").Value = 10.7
$ws4.Range("E4").Value = 16.4
$ws4.Range("F4").Value = 5.66
$ws4.Range("G4").Value = 6.22
$ws4.Range("H4").Value = 3.63

$ws4.Range("B5").Value = "RL"
$ws4.Range("C5").Value = 21.8
$ws4.Range("D5").Value = 3.7
$ws4.Range("E5").Value = 4.92
$ws4.Range("F5").Value = 9.12
$ws4.Range("G5").Value = 9.57
$ws4.Range("H5").Value = 9.91

$ws4.Range("B6").Value = "WL"
$ws4.Range("C6").Value = 18.6
$ws4.Range("D6").Value = 13.9
$ws4.Range("E6").Value = 4.01
$ws4.Range("F6").Value = 15.9
$ws4.Range("G6").Value = 11.3
$ws4.Range("H6").Value = 18.1

$ws4.Range("B7").Value = "RE"
$ws4.Range("C7").Value = 0.276
$ws4.Range("D7").Value = 0.225
$ws4.Range("E7").Value = 0.316
$ws4.Range("F7").Value = 0.105
$ws4.Range("G7").Value = 0.139
$ws4.Range("H7").Value = 0.279

$ws4.Range("B8").Value = "WE"
$ws4.Range("C8").Value = 0.293
$ws4.Range("D8").Value = 0.322
$ws4.Range("E8").Value = 0.309
$ws4.Range("F8").Value = 0.193
$ws4.Range("G8").Value = 0.131
$ws4.Range("H8").Value = 0.281

$ws4.Range("B9").Value = "LP"
$ws4.Range("C9").Value = 1.01
$ws4.Range("D9").Value = 3.53
$ws4.Range("E9").Value = 4.98
$ws4.Range("F9").Value = 1.85
$ws4.Range("G9").Value = 1.92
$ws4.Range("H9").Value = 0.78

[void]$ws4.Range("B2").Select()
